$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 8 (苏黎世意外自选 / 意外险): 保障额度 100万 -> 49万, 保费 892 -> 641
$ws.Range("F8").Value = "49万"
$ws.Range("G8").Value = 641

# Row 10 (意外津贴险): 保障额度 100 -> 200
$ws.Range("F10").Value = 200

# Recalculate the workbook so the total in G12 reflects the new values
$excel.Calculate()

# Update the active cell selection shown when the sheet was last saved
$ws.Range("H14").Select()
